# "Implemented the merge email" - replace the old Eren/Mikasa/Armin sample rows
# with a single mail-merge sample row (attachment filename, an "Enter email"
# placeholder that carries a real mailto hyperlink, and a name placeholder).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the merge-field sample data (written C, B, A to match the order
# the strings were originally entered in).
$ws.Range("C2").Value = "SamplePDF.pdf"
$ws.Range("B2").Value = "Enter email"
$ws.Range("A2").Value = "Foo"

# The old rows 3 and 4 (Mikasa / Armin) are no longer needed.
$ws.Range("A3:C4").ClearContents()

# Turn B2 into a clickable mailto: link (this also applies Excel's built-in
# "Hyperlink" style to the cell), then restore its display text.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:andersonpech@gmail.com", [ref]"", [ref]"", "andersonpech@gmail.com")
$ws.Range("B2").Value = "Enter email"

# Leave the selection where the author last left it.
$ws.Range("F8").Select()
